$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "cdm-no_metals_control"
$ws2.Name = "cdm-EX_zn2_e1-no_metals_zinc"

# Sheet 1 data (cdm-no_metals_control)
$ws1.Cells.Item(2,1).Value = [double]'76'
$ws1.Cells.Item(2,2).Value = 'S_Fatty_Acid_Metabolism'
$ws1.Cells.Item(2,3).Value = [double]'8.533353317702973e-15'
$ws1.Cells.Item(2,4).Value = [double]'0.5047619047619047'
$ws1.Cells.Item(2,5).Value = [double]'210'
$ws1.Cells.Item(3,1).Value = [double]'46'
$ws1.Cells.Item(3,2).Value = 'S_Transport__solvent_extrusion'
$ws1.Cells.Item(3,3).Value = [double]'5.122499672602734e-14'
$ws1.Cells.Item(3,4).Value = [double]'0.9615384615384616'
$ws1.Cells.Item(3,5).Value = [double]'26'
$ws1.Cells.Item(4,1).Value = [double]'7'
$ws1.Cells.Item(4,2).Value = 'S_Alternate_Carbon_and_Nitrogen_source__Dipeptide_Metabolism'
$ws1.Cells.Item(4,3).Value = [double]'5.07938356481965e-10'
$ws1.Cells.Item(4,4).Value = [double]'0.875'
$ws1.Cells.Item(4,5).Value = [double]'24'
$ws1.Cells.Item(5,1).Value = [double]'42'
$ws1.Cells.Item(5,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Queuosine_biosynthesis'
$ws1.Cells.Item(5,3).Value = [double]'2.820765231890751e-08'
$ws1.Cells.Item(5,4).Value = [double]'1'
$ws1.Cells.Item(5,5).Value = [double]'13'
$ws1.Cells.Item(6,1).Value = [double]'18'
$ws1.Cells.Item(6,2).Value = 'S_Aromatic_Compounds_Degradation__B_Ketoadipate_pathway'
$ws1.Cells.Item(6,3).Value = [double]'9.564882700522192e-08'
$ws1.Cells.Item(6,4).Value = [double]'0.7692307692307693'
$ws1.Cells.Item(6,5).Value = [double]'26'
$ws1.Cells.Item(7,1).Value = [double]'88'
$ws1.Cells.Item(7,2).Value = 'S_tRNA_Charging'
$ws1.Cells.Item(7,3).Value = [double]'2.427737807909724e-06'
$ws1.Cells.Item(7,4).Value = [double]'0.7391304347826086'
$ws1.Cells.Item(7,5).Value = [double]'23'
$ws1.Cells.Item(8,1).Value = [double]'77'
$ws1.Cells.Item(8,2).Value = 'S_Alternate_Carbon__Ascorbate_and_Aldarate_Metabolism'
$ws1.Cells.Item(8,3).Value = [double]'6.037500959134739e-06'
$ws1.Cells.Item(8,4).Value = [double]'1'
$ws1.Cells.Item(8,5).Value = [double]'9'
$ws1.Cells.Item(9,1).Value = [double]'36'
$ws1.Cells.Item(9,2).Value = 'S_Alternate_Carbon__Levulinate_Metabolism'
$ws1.Cells.Item(9,3).Value = [double]'2.30372095421102e-05'
$ws1.Cells.Item(9,4).Value = [double]'1'
$ws1.Cells.Item(9,5).Value = [double]'8'
$ws1.Cells.Item(10,1).Value = [double]'28'
$ws1.Cells.Item(10,2).Value = 'S_Glycerophospholipid_Metabolism'
$ws1.Cells.Item(10,3).Value = [double]'0.0001037616147834067'
$ws1.Cells.Item(10,4).Value = [double]'0.4'
$ws1.Cells.Item(10,5).Value = [double]'155'
$ws1.Cells.Item(11,1).Value = [double]'33'
$ws1.Cells.Item(11,2).Value = 'S_Aromatic_Compounds_Degradation__Phenylacetyl_CoA_Catabolom'
$ws1.Cells.Item(11,3).Value = [double]'0.0002006324750423651'
$ws1.Cells.Item(11,4).Value = [double]'0.7692307692307693'
$ws1.Cells.Item(11,5).Value = [double]'13'
$ws1.Cells.Item(12,1).Value = [double]'13'
$ws1.Cells.Item(12,2).Value = 'S_Alginate_biosynthesis'
$ws1.Cells.Item(12,3).Value = [double]'0.0003613274156951727'
$ws1.Cells.Item(12,4).Value = [double]'0.525'
$ws1.Cells.Item(12,5).Value = [double]'40'
$ws1.Cells.Item(13,1).Value = [double]'19'
$ws1.Cells.Item(13,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Porphyrin_and_Chlorophyll_Metabolism'
$ws1.Cells.Item(13,3).Value = [double]'0.0004778202916086266'
$ws1.Cells.Item(13,4).Value = [double]'0.5263157894736842'
$ws1.Cells.Item(13,5).Value = [double]'38'
$ws1.Cells.Item(14,1).Value = [double]'60'
$ws1.Cells.Item(14,2).Value = 'Unassigned'
$ws1.Cells.Item(14,3).Value = [double]'0.000846117971279241'
$ws1.Cells.Item(14,4).Value = [double]'0.5'
$ws1.Cells.Item(14,5).Value = [double]'42'
$ws1.Cells.Item(15,1).Value = [double]'24'
$ws1.Cells.Item(15,2).Value = 'S_Alternate_Carbon_and_Nitrogen_source__D_Amino_acids_Metabolism'
$ws1.Cells.Item(15,3).Value = [double]'0.001230422132652071'
$ws1.Cells.Item(15,4).Value = [double]'0.6666666666666666'
$ws1.Cells.Item(15,5).Value = [double]'15'
$ws1.Cells.Item(16,1).Value = [double]'83'
$ws1.Cells.Item(16,2).Value = 'Aliphatic open-chain ketones metabolism'
$ws1.Cells.Item(16,3).Value = [double]'0.001272445350034012'
$ws1.Cells.Item(16,4).Value = [double]'1'
$ws1.Cells.Item(16,5).Value = [double]'5'
$ws1.Cells.Item(17,1).Value = [double]'5'
$ws1.Cells.Item(17,2).Value = 'S_Heavy_Metal_Tolerance'
$ws1.Cells.Item(17,3).Value = [double]'0.003494952742827158'
$ws1.Cells.Item(17,4).Value = [double]'0.5'
$ws1.Cells.Item(17,5).Value = [double]'32'
$ws1.Cells.Item(18,1).Value = [double]'2'
$ws1.Cells.Item(18,2).Value = 'S_Aromatic_Compounds_Degradation__Gallic_acid_Metabolism'
$ws1.Cells.Item(18,3).Value = [double]'0.004836616070415369'
$ws1.Cells.Item(18,4).Value = [double]'1'
$ws1.Cells.Item(18,5).Value = [double]'4'
$ws1.Cells.Item(19,1).Value = [double]'80'
$ws1.Cells.Item(19,2).Value = 'S_Plant_growth_promoting'
$ws1.Cells.Item(19,3).Value = [double]'0.004836616070415369'
$ws1.Cells.Item(19,4).Value = [double]'1'
$ws1.Cells.Item(19,5).Value = [double]'4'
$ws1.Cells.Item(20,1).Value = [double]'38'
$ws1.Cells.Item(20,2).Value = 'S_Nitrogen_Metabolism'
$ws1.Cells.Item(20,3).Value = [double]'0.004836616070415369'
$ws1.Cells.Item(20,4).Value = [double]'1'
$ws1.Cells.Item(20,5).Value = [double]'4'
$ws1.Cells.Item(21,1).Value = [double]'40'
$ws1.Cells.Item(21,2).Value = 'S_Aromatic_Compounds_Degradation__Nicotinic_acid_pathway'
$ws1.Cells.Item(21,3).Value = [double]'0.00596246465868094'
$ws1.Cells.Item(21,4).Value = [double]'0.8333333333333334'
$ws1.Cells.Item(21,5).Value = [double]'6'
$ws1.Cells.Item(22,1).Value = [double]'3'
$ws1.Cells.Item(22,2).Value = 'S_Alternate_Carbon'
$ws1.Cells.Item(22,3).Value = [double]'0.007960890625810594'
$ws1.Cells.Item(22,4).Value = [double]'0.6153846153846154'
$ws1.Cells.Item(22,5).Value = [double]'13'
$ws1.Cells.Item(23,1).Value = [double]'67'
$ws1.Cells.Item(23,2).Value = 'S_Histidine_Metabolism'
$ws1.Cells.Item(23,3).Value = [double]'0.01756593556708986'
$ws1.Cells.Item(23,4).Value = [double]'0.5294117647058824'
$ws1.Cells.Item(23,5).Value = [double]'17'
$ws1.Cells.Item(24,1).Value = [double]'78'
$ws1.Cells.Item(24,2).Value = 'S_Alternate_Carbon_and_Nitrogen_source__Hydroxy_proline_metabolism'
$ws1.Cells.Item(24,3).Value = [double]'0.01836657842843447'
$ws1.Cells.Item(24,4).Value = [double]'1'
$ws1.Cells.Item(24,5).Value = [double]'3'
$ws1.Cells.Item(25,1).Value = [double]'35'
$ws1.Cells.Item(25,2).Value = 'S_Sulfur_Metabolism'
$ws1.Cells.Item(25,3).Value = [double]'0.02190214730732861'
$ws1.Cells.Item(25,4).Value = [double]'0.4782608695652174'
$ws1.Cells.Item(25,5).Value = [double]'23'
$ws1.Cells.Item(26,1).Value = [double]'50'
$ws1.Cells.Item(26,2).Value = 'S_Xenobiotic_tolerance'
$ws1.Cells.Item(26,3).Value = [double]'0.04394732399713096'
$ws1.Cells.Item(26,4).Value = [double]'0.5454545454545454'
$ws1.Cells.Item(26,5).Value = [double]'11'
$ws1.Cells.Item(27,1).Value = [double]'74'
$ws1.Cells.Item(27,2).Value = 'S_Pyruvate_Metabolism'
$ws1.Cells.Item(27,3).Value = [double]'0.04535496753846053'
$ws1.Cells.Item(27,4).Value = [double]'0.6666666666666666'
$ws1.Cells.Item(27,5).Value = [double]'6'
$ws1.Cells.Item(28,1).Value = [double]'64'
$ws1.Cells.Item(28,2).Value = 'S_Nucleotide_Salvage_Pathway'
$ws1.Cells.Item(28,3).Value = [double]'0.04535496753846053'
$ws1.Cells.Item(28,4).Value = [double]'0.6666666666666666'
$ws1.Cells.Item(28,5).Value = [double]'6'
$ws1.Cells.Item(29,1).Value = [double]'51'
$ws1.Cells.Item(29,2).Value = 'S_Phosphonate_and_phosphinate_metabolism'
$ws1.Cells.Item(29,3).Value = [double]'0.05895646550249176'
$ws1.Cells.Item(29,4).Value = [double]'0.75'
$ws1.Cells.Item(29,5).Value = [double]'4'
$ws1.Cells.Item(30,1).Value = [double]'0'
$ws1.Cells.Item(30,2).Value = 'Murein Recycling'
$ws1.Cells.Item(30,3).Value = [double]'0.06967865357090897'
$ws1.Cells.Item(30,4).Value = [double]'1'
$ws1.Cells.Item(30,5).Value = [double]'2'
$ws1.Cells.Item(31,1).Value = [double]'12'
$ws1.Cells.Item(31,2).Value = 'S_Alternate_Carbon_and_Nitrogen_source__Nucleotide_Metabolism'
$ws1.Cells.Item(31,3).Value = [double]'0.07310822657446485'
$ws1.Cells.Item(31,4).Value = [double]'0.4666666666666667'
$ws1.Cells.Item(31,5).Value = [double]'15'
$ws1.Cells.Item(32,1).Value = [double]'4'
$ws1.Cells.Item(32,2).Value = 'S_Cell_Envelope_Biosynthesis__Peptidoglycan_Biosynthesis'
$ws1.Cells.Item(32,3).Value = [double]'0.10047257484093'
$ws1.Cells.Item(32,4).Value = [double]'0.3529411764705883'
$ws1.Cells.Item(32,5).Value = [double]'51'
$ws1.Cells.Item(33,1).Value = [double]'23'
$ws1.Cells.Item(33,2).Value = 'S_Methionine_Metabolism'
$ws1.Cells.Item(33,3).Value = [double]'0.1375833926400542'
$ws1.Cells.Item(33,4).Value = [double]'0.4285714285714285'
$ws1.Cells.Item(33,5).Value = [double]'14'
$ws1.Cells.Item(34,1).Value = [double]'69'
$ws1.Cells.Item(34,2).Value = 'S_Cell_Envelope_Biosynthesis__Cellulose_Metabolism'
$ws1.Cells.Item(34,3).Value = [double]'0.172302803855858'
$ws1.Cells.Item(34,4).Value = [double]'0.6666666666666666'
$ws1.Cells.Item(34,5).Value = [double]'3'
$ws1.Cells.Item(35,1).Value = [double]'85'
$ws1.Cells.Item(35,2).Value = 'S_Alternate_Carbon__Fructose_Metabolism'
$ws1.Cells.Item(35,3).Value = [double]'0.172302803855858'
$ws1.Cells.Item(35,4).Value = [double]'0.6666666666666666'
$ws1.Cells.Item(35,5).Value = [double]'3'
$ws1.Cells.Item(36,1).Value = [double]'87'
$ws1.Cells.Item(36,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Vitamin_B6_Metabolism'
$ws1.Cells.Item(36,3).Value = [double]'0.2576635754315496'
$ws1.Cells.Item(36,4).Value = [double]'0.4'
$ws1.Cells.Item(36,5).Value = [double]'10'
$ws1.Cells.Item(37,1).Value = [double]'79'
$ws1.Cells.Item(37,2).Value = 'S_Phenylalanine_Tyrosine_Tryptophan_Biosynthesis'
$ws1.Cells.Item(37,3).Value = [double]'0.2671151394417686'
$ws1.Cells.Item(37,4).Value = [double]'0.3333333333333333'
$ws1.Cells.Item(37,5).Value = [double]'27'
$ws1.Cells.Item(38,1).Value = [double]'81'
$ws1.Cells.Item(38,2).Value = 'S_Alternate_Carbon_and_Nitrogen_source__Ectoine_Metabolism'
$ws1.Cells.Item(38,3).Value = [double]'0.2856491422092242'
$ws1.Cells.Item(38,4).Value = [double]'0.5'
$ws1.Cells.Item(38,5).Value = [double]'4'
$ws1.Cells.Item(39,1).Value = [double]'25'
$ws1.Cells.Item(39,2).Value = 'S_Alternate_Carbon_and_Nitrogen_source__Amines_and_Polyamines_Metabolism'
$ws1.Cells.Item(39,3).Value = [double]'0.3586202445561254'
$ws1.Cells.Item(39,4).Value = [double]'0.3'
$ws1.Cells.Item(39,5).Value = [double]'40'
$ws1.Cells.Item(40,1).Value = [double]'22'
$ws1.Cells.Item(40,2).Value = 'S_Alternate_Carbon__Butanediol_Metabolism'
$ws1.Cells.Item(40,3).Value = [double]'0.4585072022541679'
$ws1.Cells.Item(40,4).Value = [double]'0.5'
$ws1.Cells.Item(40,5).Value = [double]'2'
$ws1.Cells.Item(41,1).Value = [double]'17'
$ws1.Cells.Item(41,2).Value = 'S_Transport__Inner_Membrane'
$ws1.Cells.Item(41,3).Value = [double]'0.4834791060654365'
$ws1.Cells.Item(41,4).Value = [double]'0.2672811059907834'
$ws1.Cells.Item(41,5).Value = [double]'217'
$ws1.Cells.Item(42,1).Value = [double]'55'
$ws1.Cells.Item(42,2).Value = 'S_Formaldehyde_Metabolism'
$ws1.Cells.Item(42,3).Value = [double]'0.4992979086356545'
$ws1.Cells.Item(42,4).Value = [double]'0.3333333333333333'
$ws1.Cells.Item(42,5).Value = [double]'6'
$ws1.Cells.Item(43,1).Value = [double]'41'
$ws1.Cells.Item(43,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Molybdopterin_Biosynthesis'
$ws1.Cells.Item(43,3).Value = [double]'0.527997820125692'
$ws1.Cells.Item(43,4).Value = [double]'0.2857142857142857'
$ws1.Cells.Item(43,5).Value = [double]'14'
$ws1.Cells.Item(44,1).Value = [double]'56'
$ws1.Cells.Item(44,2).Value = 'S_Pyrimidine_Metabolism'
$ws1.Cells.Item(44,3).Value = [double]'0.5807867915859364'
$ws1.Cells.Item(44,4).Value = [double]'0.26'
$ws1.Cells.Item(44,5).Value = [double]'50'
$ws1.Cells.Item(45,1).Value = [double]'20'
$ws1.Cells.Item(45,2).Value = 'S_Starch_and_Sucrose_Metabolism'
$ws1.Cells.Item(45,3).Value = [double]'0.5888537690237506'
$ws1.Cells.Item(45,4).Value = [double]'0.2666666666666667'
$ws1.Cells.Item(45,5).Value = [double]'15'
$ws1.Cells.Item(46,1).Value = [double]'15'
$ws1.Cells.Item(46,2).Value = 'S_Valine__Leucine__and_Isoleucine_Metabolism'
$ws1.Cells.Item(46,3).Value = [double]'0.6018438180271301'
$ws1.Cells.Item(46,4).Value = [double]'0.2571428571428571'
$ws1.Cells.Item(46,5).Value = [double]'35'
$ws1.Cells.Item(47,1).Value = [double]'61'
$ws1.Cells.Item(47,2).Value = 'S_Cysteine_Metabolism'
$ws1.Cells.Item(47,3).Value = [double]'0.6673749332079857'
$ws1.Cells.Item(47,4).Value = [double]'0.25'
$ws1.Cells.Item(47,5).Value = [double]'8'
$ws1.Cells.Item(48,1).Value = [double]'59'
$ws1.Cells.Item(48,2).Value = 'S_Alanine_and_Aspartate_Metabolism'
$ws1.Cells.Item(48,3).Value = [double]'0.6673749332079857'
$ws1.Cells.Item(48,4).Value = [double]'0.25'
$ws1.Cells.Item(48,5).Value = [double]'8'
$ws1.Cells.Item(49,1).Value = [double]'21'
$ws1.Cells.Item(49,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis'
$ws1.Cells.Item(49,3).Value = [double]'0.7069294878680226'
$ws1.Cells.Item(49,4).Value = [double]'0.25'
$ws1.Cells.Item(49,5).Value = [double]'4'
$ws1.Cells.Item(50,1).Value = [double]'62'
$ws1.Cells.Item(50,2).Value = 'S_Transport__Outer_Membrane'
$ws1.Cells.Item(50,3).Value = [double]'0.7193323514820571'
$ws1.Cells.Item(50,4).Value = [double]'0.2516556291390729'
$ws1.Cells.Item(50,5).Value = [double]'302'
$ws1.Cells.Item(51,1).Value = [double]'27'
$ws1.Cells.Item(51,2).Value = 'S_Glycine__Serine_and_threonine_metabolism'
$ws1.Cells.Item(51,3).Value = [double]'0.8005652925441062'
$ws1.Cells.Item(51,4).Value = [double]'0.2'
$ws1.Cells.Item(51,5).Value = [double]'15'
$ws1.Cells.Item(52,1).Value = [double]'32'
$ws1.Cells.Item(52,2).Value = 'S_Cell_Envelope_Biosynthesis__O_antigen_Biosynthesis'
$ws1.Cells.Item(52,3).Value = [double]'0.8005652925441062'
$ws1.Cells.Item(52,4).Value = [double]'0.2'
$ws1.Cells.Item(52,5).Value = [double]'15'
$ws1.Cells.Item(53,1).Value = [double]'8'
$ws1.Cells.Item(53,2).Value = 'S_Butanoate_Metabolism'
$ws1.Cells.Item(53,3).Value = [double]'0.8414603220359542'
$ws1.Cells.Item(53,4).Value = [double]'0.1666666666666667'
$ws1.Cells.Item(53,5).Value = [double]'6'
$ws1.Cells.Item(54,1).Value = [double]'72'
$ws1.Cells.Item(54,2).Value = 'S_Cell_Envelope_Biosynthesis'
$ws1.Cells.Item(54,3).Value = [double]'0.8953470941245922'
$ws1.Cells.Item(54,4).Value = [double]'0.1538461538461539'
$ws1.Cells.Item(54,5).Value = [double]'13'
$ws1.Cells.Item(55,1).Value = [double]'71'
$ws1.Cells.Item(55,2).Value = 'S_Lysine_Metabolism'
$ws1.Cells.Item(55,3).Value = [double]'0.9036123201281688'
$ws1.Cells.Item(55,4).Value = [double]'0.1818181818181818'
$ws1.Cells.Item(55,5).Value = [double]'33'
$ws1.Cells.Item(56,1).Value = [double]'65'
$ws1.Cells.Item(56,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__FeS_clusters_metabolism'
$ws1.Cells.Item(56,3).Value = [double]'0.9142784757649196'
$ws1.Cells.Item(56,4).Value = [double]'0.125'
$ws1.Cells.Item(56,5).Value = [double]'8'
$ws1.Cells.Item(57,1).Value = [double]'1'
$ws1.Cells.Item(57,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Thiamine_Metabolism'
$ws1.Cells.Item(57,3).Value = [double]'0.9142784757649196'
$ws1.Cells.Item(57,4).Value = [double]'0.125'
$ws1.Cells.Item(57,5).Value = [double]'8'
$ws1.Cells.Item(58,1).Value = [double]'49'
$ws1.Cells.Item(58,2).Value = 'S_Cell_Envelope_Biosynthesis__LPS_Biosynthesis'
$ws1.Cells.Item(58,3).Value = [double]'0.918201111838685'
$ws1.Cells.Item(58,4).Value = [double]'0.1428571428571428'
$ws1.Cells.Item(58,5).Value = [double]'14'
$ws1.Cells.Item(59,1).Value = [double]'82'
$ws1.Cells.Item(59,2).Value = 'S_Cell_Envelope_Biosynthesis__Lipid_A_Biosynthesis'
$ws1.Cells.Item(59,3).Value = [double]'0.9362989301916973'
$ws1.Cells.Item(59,4).Value = [double]'0.1333333333333333'
$ws1.Cells.Item(59,5).Value = [double]'15'
$ws1.Cells.Item(60,1).Value = [double]'84'
$ws1.Cells.Item(60,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Ubiquinone_biosynthesis'
$ws1.Cells.Item(60,3).Value = [double]'0.9536737105512071'
$ws1.Cells.Item(60,4).Value = [double]'0.1'
$ws1.Cells.Item(60,5).Value = [double]'10'
$ws1.Cells.Item(61,1).Value = [double]'89'
$ws1.Cells.Item(61,2).Value = 'S_Transport__ABC_system'
$ws1.Cells.Item(61,3).Value = [double]'0.9604361727967153'
$ws1.Cells.Item(61,4).Value = [double]'0.2043795620437956'
$ws1.Cells.Item(61,5).Value = [double]'137'
$ws1.Cells.Item(62,1).Value = [double]'29'
$ws1.Cells.Item(62,2).Value = 'S_Arginine_and_Proline_Metabolism'
$ws1.Cells.Item(62,3).Value = [double]'0.9612669257067126'
$ws1.Cells.Item(62,4).Value = [double]'0.1428571428571428'
$ws1.Cells.Item(62,5).Value = [double]'28'
$ws1.Cells.Item(63,1).Value = [double]'9'
$ws1.Cells.Item(63,2).Value = 'S_Glyoxylate_and_dicarboxylate_metabolism'
$ws1.Cells.Item(63,3).Value = [double]'0.9659500978477161'
$ws1.Cells.Item(63,4).Value = [double]'0.09090909090909091'
$ws1.Cells.Item(63,5).Value = [double]'11'
$ws1.Cells.Item(64,1).Value = [double]'44'
$ws1.Cells.Item(64,2).Value = 'S_Purine_Metabolism'
$ws1.Cells.Item(64,3).Value = [double]'0.9674495471864356'
$ws1.Cells.Item(64,4).Value = [double]'0.1818181818181818'
$ws1.Cells.Item(64,5).Value = [double]'77'
$ws1.Cells.Item(65,1).Value = [double]'68'
$ws1.Cells.Item(65,2).Value = 'S_Glycolysis'
$ws1.Cells.Item(65,3).Value = [double]'0.9960678662826885'
$ws1.Cells.Item(65,4).Value = [double]'0.05555555555555555'
$ws1.Cells.Item(65,5).Value = [double]'18'
$ws1.Cells.Item(66,1).Value = [double]'75'
$ws1.Cells.Item(66,2).Value = 'S_Aromatic_Compounds_Degradation__Toluene_Pathway'
$ws1.Cells.Item(66,3).Value = [double]'0.9988588030145974'
$ws1.Cells.Item(66,4).Value = [double]'0.06666666666666667'
$ws1.Cells.Item(66,5).Value = [double]'30'
$ws1.Cells.Item(67,1).Value = [double]'34'
$ws1.Cells.Item(67,2).Value = 'S_Fatty_Acid__Biosynthesis'
$ws1.Cells.Item(67,3).Value = [double]'0.999815915265249'
$ws1.Cells.Item(67,4).Value = [double]'0.1066666666666667'
$ws1.Cells.Item(67,5).Value = [double]'75'
$ws1.Cells.Item(68,1).Value = [double]'57'
$ws1.Cells.Item(68,2).Value = 'S_Iron_uptake_and_metabolism'
$ws1.Cells.Item(68,3).Value = [double]'0.9999997874161561'
$ws1.Cells.Item(68,4).Value = [double]'0.05263157894736842'
$ws1.Cells.Item(68,5).Value = [double]'76'
$ws1.Cells.Item(69,1).Value = [double]'70'
$ws1.Cells.Item(69,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Pantothenate_and_CoA_Biosynthesis'
$ws1.Cells.Item(69,3).Value = [double]'1'
$ws1.Cells.Item(69,4).Value = [double]'0'
$ws1.Cells.Item(69,5).Value = [double]'10'
$ws1.Cells.Item(70,1).Value = [double]'73'
$ws1.Cells.Item(70,2).Value = 'S_Inorganic_polyphosphates_metabolism'
$ws1.Cells.Item(70,3).Value = [double]'1'
$ws1.Cells.Item(70,4).Value = [double]'0'
$ws1.Cells.Item(70,5).Value = [double]'7'
$ws1.Cells.Item(71,1).Value = [double]'86'
$ws1.Cells.Item(71,2).Value = 'S_TCA_Cycle'
$ws1.Cells.Item(71,3).Value = [double]'1'
$ws1.Cells.Item(71,4).Value = [double]'0'
$ws1.Cells.Item(71,5).Value = [double]'20'
$ws1.Cells.Item(72,1).Value = [double]'45'
$ws1.Cells.Item(72,2).Value = 'S_Alternate_Carbon__Propanoate_Metabolism'
$ws1.Cells.Item(72,3).Value = [double]'1'
$ws1.Cells.Item(72,4).Value = [double]'0'
$ws1.Cells.Item(72,5).Value = [double]'7'
$ws1.Cells.Item(73,1).Value = [double]'63'
$ws1.Cells.Item(73,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Riboflavin_Metabolism'
$ws1.Cells.Item(73,3).Value = [double]'1'
$ws1.Cells.Item(73,4).Value = [double]'0'
$ws1.Cells.Item(73,5).Value = [double]'9'
$ws1.Cells.Item(74,1).Value = [double]'6'
$ws1.Cells.Item(74,2).Value = 'Extracellular exchange'
$ws1.Cells.Item(74,3).Value = [double]'1'
$ws1.Cells.Item(74,4).Value = [double]'0'
$ws1.Cells.Item(74,5).Value = [double]'348'
$ws1.Cells.Item(75,1).Value = [double]'10'
$ws1.Cells.Item(75,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Nicotinamide_Biosynthesis'
$ws1.Cells.Item(75,3).Value = [double]'1'
$ws1.Cells.Item(75,4).Value = [double]'0'
$ws1.Cells.Item(75,5).Value = [double]'13'
$ws1.Cells.Item(76,1).Value = [double]'11'
$ws1.Cells.Item(76,2).Value = 'S_Cell_Envelope_Biosynthesis__Biosynthesis_of_L_glycero_D_manno_heptose__Hep_'
$ws1.Cells.Item(76,3).Value = [double]'1'
$ws1.Cells.Item(76,4).Value = [double]'0'
$ws1.Cells.Item(76,5).Value = [double]'5'
$ws1.Cells.Item(77,1).Value = [double]'14'
$ws1.Cells.Item(77,2).Value = 'S_Oxidative_Phosphorylation'
$ws1.Cells.Item(77,3).Value = [double]'1'
$ws1.Cells.Item(77,4).Value = [double]'0'
$ws1.Cells.Item(77,5).Value = [double]'10'
$ws1.Cells.Item(78,1).Value = [double]'16'
$ws1.Cells.Item(78,2).Value = 'Intracellular source/sink'
$ws1.Cells.Item(78,3).Value = [double]'1'
$ws1.Cells.Item(78,4).Value = [double]'0'
$ws1.Cells.Item(78,5).Value = [double]'2'
$ws1.Cells.Item(79,1).Value = [double]'26'
$ws1.Cells.Item(79,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Folate_Biosynthesis'
$ws1.Cells.Item(79,3).Value = [double]'1'
$ws1.Cells.Item(79,4).Value = [double]'0'
$ws1.Cells.Item(79,5).Value = [double]'14'
$ws1.Cells.Item(80,1).Value = [double]'30'
$ws1.Cells.Item(80,2).Value = 'S_Alternate_Carbon__Ribose_Metabolism'
$ws1.Cells.Item(80,3).Value = [double]'1'
$ws1.Cells.Item(80,4).Value = [double]'0'
$ws1.Cells.Item(80,5).Value = [double]'5'
$ws1.Cells.Item(81,1).Value = [double]'31'
$ws1.Cells.Item(81,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Biotin_Biosynthesis'
$ws1.Cells.Item(81,3).Value = [double]'1'
$ws1.Cells.Item(81,4).Value = [double]'0'
$ws1.Cells.Item(81,5).Value = [double]'15'
$ws1.Cells.Item(82,1).Value = [double]'37'
$ws1.Cells.Item(82,2).Value = 'S_Gluconeogenesis'
$ws1.Cells.Item(82,3).Value = [double]'1'
$ws1.Cells.Item(82,4).Value = [double]'0'
$ws1.Cells.Item(82,5).Value = [double]'12'
$ws1.Cells.Item(83,1).Value = [double]'39'
$ws1.Cells.Item(83,2).Value = 'S_Pentose_Phosphate_Pathway'
$ws1.Cells.Item(83,3).Value = [double]'1'
$ws1.Cells.Item(83,4).Value = [double]'0'
$ws1.Cells.Item(83,5).Value = [double]'7'
$ws1.Cells.Item(84,1).Value = [double]'43'
$ws1.Cells.Item(84,2).Value = 'S_PHAs_Metabolism'
$ws1.Cells.Item(84,3).Value = [double]'1'
$ws1.Cells.Item(84,4).Value = [double]'0.02631578947368421'
$ws1.Cells.Item(84,5).Value = [double]'152'
$ws1.Cells.Item(85,1).Value = [double]'47'
$ws1.Cells.Item(85,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis___Pyrroloquinoline_quinone_biosynthesis'
$ws1.Cells.Item(85,3).Value = [double]'1'
$ws1.Cells.Item(85,4).Value = [double]'0'
$ws1.Cells.Item(85,5).Value = [double]'6'
$ws1.Cells.Item(86,1).Value = [double]'48'
$ws1.Cells.Item(86,2).Value = 'S_Aromatic_Compounds_Degradation__Homogentisate_pathway'
$ws1.Cells.Item(86,3).Value = [double]'1'
$ws1.Cells.Item(86,4).Value = [double]'0'
$ws1.Cells.Item(86,5).Value = [double]'6'
$ws1.Cells.Item(87,1).Value = [double]'52'
$ws1.Cells.Item(87,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Terpenoid_backbone_biosynthesis'
$ws1.Cells.Item(87,3).Value = [double]'1'
$ws1.Cells.Item(87,4).Value = [double]'0'
$ws1.Cells.Item(87,5).Value = [double]'14'
$ws1.Cells.Item(88,1).Value = [double]'53'
$ws1.Cells.Item(88,2).Value = 'Intracellular demand'
$ws1.Cells.Item(88,3).Value = [double]'1'
$ws1.Cells.Item(88,4).Value = [double]'0'
$ws1.Cells.Item(88,5).Value = [double]'31'
$ws1.Cells.Item(89,1).Value = [double]'54'
$ws1.Cells.Item(89,2).Value = 'Biomass and maintenance functions'
$ws1.Cells.Item(89,3).Value = [double]'1'
$ws1.Cells.Item(89,4).Value = [double]'0'
$ws1.Cells.Item(89,5).Value = [double]'4'
$ws1.Cells.Item(90,1).Value = [double]'58'
$ws1.Cells.Item(90,2).Value = 'S_Urea_cycleamino_group_metabolism'
$ws1.Cells.Item(90,3).Value = [double]'1'
$ws1.Cells.Item(90,4).Value = [double]'0'
$ws1.Cells.Item(90,5).Value = [double]'7'
$ws1.Cells.Item(91,1).Value = [double]'66'
$ws1.Cells.Item(91,2).Value = 'S_Glutamate_Metabolism'
$ws1.Cells.Item(91,3).Value = [double]'1'
$ws1.Cells.Item(91,4).Value = [double]'0'
$ws1.Cells.Item(91,5).Value = [double]'17'
$ws1.Cells.Item(92,1).Value = [double]'90'
$ws1.Cells.Item(92,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__One_Carbon_pool_by_folate'
$ws1.Cells.Item(92,3).Value = [double]'1'
$ws1.Cells.Item(92,4).Value = [double]'0'
$ws1.Cells.Item(92,5).Value = [double]'7'

# Sheet 2 data (cdm-EX_zn2_e1-no_metals_zinc)
$ws2.Cells.Item(2,1).Value = [double]'76'
$ws2.Cells.Item(2,2).Value = 'S_Fatty_Acid_Metabolism'
$ws2.Cells.Item(2,3).Value = [double]'6.779324403041417e-18'
$ws2.Cells.Item(2,4).Value = [double]'0.4666666666666667'
$ws2.Cells.Item(2,5).Value = [double]'210'
$ws2.Cells.Item(3,1).Value = [double]'89'
$ws2.Cells.Item(3,2).Value = 'S_Transport__ABC_system'
$ws2.Cells.Item(3,3).Value = [double]'1.231585793241232e-13'
$ws2.Cells.Item(3,4).Value = [double]'0.489051094890511'
$ws2.Cells.Item(3,5).Value = [double]'137'
$ws2.Cells.Item(4,1).Value = [double]'46'
$ws2.Cells.Item(4,2).Value = 'S_Transport__solvent_extrusion'
$ws2.Cells.Item(4,3).Value = [double]'1.138385645187223e-10'
$ws2.Cells.Item(4,4).Value = [double]'0.8076923076923077'
$ws2.Cells.Item(4,5).Value = [double]'26'
$ws2.Cells.Item(5,1).Value = [double]'25'
$ws2.Cells.Item(5,2).Value = 'S_Alternate_Carbon_and_Nitrogen_source__Amines_and_Polyamines_Metabolism'
$ws2.Cells.Item(5,3).Value = [double]'1.074112203646232e-05'
$ws2.Cells.Item(5,4).Value = [double]'0.525'
$ws2.Cells.Item(5,5).Value = [double]'40'
$ws2.Cells.Item(6,1).Value = [double]'13'
$ws2.Cells.Item(6,2).Value = 'S_Alginate_biosynthesis'
$ws2.Cells.Item(6,3).Value = [double]'1.074112203646232e-05'
$ws2.Cells.Item(6,4).Value = [double]'0.525'
$ws2.Cells.Item(6,5).Value = [double]'40'
$ws2.Cells.Item(7,1).Value = [double]'77'
$ws2.Cells.Item(7,2).Value = 'S_Alternate_Carbon__Ascorbate_and_Aldarate_Metabolism'
$ws2.Cells.Item(7,3).Value = [double]'2.827555945472601e-05'
$ws2.Cells.Item(7,4).Value = [double]'0.8888888888888888'
$ws2.Cells.Item(7,5).Value = [double]'9'
$ws2.Cells.Item(8,1).Value = [double]'19'
$ws2.Cells.Item(8,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Porphyrin_and_Chlorophyll_Metabolism'
$ws2.Cells.Item(8,3).Value = [double]'6.905719648533999e-05'
$ws2.Cells.Item(8,4).Value = [double]'0.5'
$ws2.Cells.Item(8,5).Value = [double]'38'
$ws2.Cells.Item(9,1).Value = [double]'24'
$ws2.Cells.Item(9,2).Value = 'S_Alternate_Carbon_and_Nitrogen_source__D_Amino_acids_Metabolism'
$ws2.Cells.Item(9,3).Value = [double]'0.0001777847557186491'
$ws2.Cells.Item(9,4).Value = [double]'0.6666666666666666'
$ws2.Cells.Item(9,5).Value = [double]'15'
$ws2.Cells.Item(10,1).Value = [double]'12'
$ws2.Cells.Item(10,2).Value = 'S_Alternate_Carbon_and_Nitrogen_source__Nucleotide_Metabolism'
$ws2.Cells.Item(10,3).Value = [double]'0.001166565680518156'
$ws2.Cells.Item(10,4).Value = [double]'0.6'
$ws2.Cells.Item(10,5).Value = [double]'15'
$ws2.Cells.Item(11,1).Value = [double]'38'
$ws2.Cells.Item(11,2).Value = 'S_Nitrogen_Metabolism'
$ws2.Cells.Item(11,3).Value = [double]'0.001984927814041523'
$ws2.Cells.Item(11,4).Value = [double]'1'
$ws2.Cells.Item(11,5).Value = [double]'4'
$ws2.Cells.Item(12,1).Value = [double]'85'
$ws2.Cells.Item(12,2).Value = 'S_Alternate_Carbon__Fructose_Metabolism'
$ws2.Cells.Item(12,3).Value = [double]'0.009421962545872421'
$ws2.Cells.Item(12,4).Value = [double]'1'
$ws2.Cells.Item(12,5).Value = [double]'3'
$ws2.Cells.Item(13,1).Value = [double]'78'
$ws2.Cells.Item(13,2).Value = 'S_Alternate_Carbon_and_Nitrogen_source__Hydroxy_proline_metabolism'
$ws2.Cells.Item(13,3).Value = [double]'0.009421962545872421'
$ws2.Cells.Item(13,4).Value = [double]'1'
$ws2.Cells.Item(13,5).Value = [double]'3'
$ws2.Cells.Item(14,1).Value = [double]'36'
$ws2.Cells.Item(14,2).Value = 'S_Alternate_Carbon__Levulinate_Metabolism'
$ws2.Cells.Item(14,3).Value = [double]'0.01317622558189333'
$ws2.Cells.Item(14,4).Value = [double]'0.625'
$ws2.Cells.Item(14,5).Value = [double]'8'
$ws2.Cells.Item(15,1).Value = [double]'48'
$ws2.Cells.Item(15,2).Value = 'S_Aromatic_Compounds_Degradation__Homogentisate_pathway'
$ws2.Cells.Item(15,3).Value = [double]'0.02062837744710942'
$ws2.Cells.Item(15,4).Value = [double]'0.6666666666666666'
$ws2.Cells.Item(15,5).Value = [double]'6'
$ws2.Cells.Item(16,1).Value = [double]'33'
$ws2.Cells.Item(16,2).Value = 'S_Aromatic_Compounds_Degradation__Phenylacetyl_CoA_Catabolom'
$ws2.Cells.Item(16,3).Value = [double]'0.03838108686849758'
$ws2.Cells.Item(16,4).Value = [double]'0.4615384615384616'
$ws2.Cells.Item(16,5).Value = [double]'13'
$ws2.Cells.Item(17,1).Value = [double]'22'
$ws2.Cells.Item(17,2).Value = 'S_Alternate_Carbon__Butanediol_Metabolism'
$ws2.Cells.Item(17,3).Value = [double]'0.04466651612103215'
$ws2.Cells.Item(17,4).Value = [double]'1'
$ws2.Cells.Item(17,5).Value = [double]'2'
$ws2.Cells.Item(18,1).Value = [double]'28'
$ws2.Cells.Item(18,2).Value = 'S_Glycerophospholipid_Metabolism'
$ws2.Cells.Item(18,3).Value = [double]'0.06192863658221925'
$ws2.Cells.Item(18,4).Value = [double]'0.2645161290322581'
$ws2.Cells.Item(18,5).Value = [double]'155'
$ws2.Cells.Item(19,1).Value = [double]'18'
$ws2.Cells.Item(19,2).Value = 'S_Aromatic_Compounds_Degradation__B_Ketoadipate_pathway'
$ws2.Cells.Item(19,3).Value = [double]'0.07897969482351258'
$ws2.Cells.Item(19,4).Value = [double]'0.3461538461538461'
$ws2.Cells.Item(19,5).Value = [double]'26'
$ws2.Cells.Item(20,1).Value = [double]'15'
$ws2.Cells.Item(20,2).Value = 'S_Valine__Leucine__and_Isoleucine_Metabolism'
$ws2.Cells.Item(20,3).Value = [double]'0.101838431156142'
$ws2.Cells.Item(20,4).Value = [double]'0.3142857142857143'
$ws2.Cells.Item(20,5).Value = [double]'35'
$ws2.Cells.Item(21,1).Value = [double]'4'
$ws2.Cells.Item(21,2).Value = 'S_Cell_Envelope_Biosynthesis__Peptidoglycan_Biosynthesis'
$ws2.Cells.Item(21,3).Value = [double]'0.1022469311817887'
$ws2.Cells.Item(21,4).Value = [double]'0.2941176470588235'
$ws2.Cells.Item(21,5).Value = [double]'51'
$ws2.Cells.Item(22,1).Value = [double]'40'
$ws2.Cells.Item(22,2).Value = 'S_Aromatic_Compounds_Degradation__Nicotinic_acid_pathway'
$ws2.Cells.Item(22,3).Value = [double]'0.1132745913130778'
$ws2.Cells.Item(22,4).Value = [double]'0.5'
$ws2.Cells.Item(22,5).Value = [double]'6'
$ws2.Cells.Item(23,1).Value = [double]'71'
$ws2.Cells.Item(23,2).Value = 'S_Lysine_Metabolism'
$ws2.Cells.Item(23,3).Value = [double]'0.1406118011665893'
$ws2.Cells.Item(23,4).Value = [double]'0.303030303030303'
$ws2.Cells.Item(23,5).Value = [double]'33'
$ws2.Cells.Item(24,1).Value = [double]'81'
$ws2.Cells.Item(24,2).Value = 'S_Alternate_Carbon_and_Nitrogen_source__Ectoine_Metabolism'
$ws2.Cells.Item(24,3).Value = [double]'0.198578179801338'
$ws2.Cells.Item(24,4).Value = [double]'0.5'
$ws2.Cells.Item(24,5).Value = [double]'4'
$ws2.Cells.Item(25,1).Value = [double]'80'
$ws2.Cells.Item(25,2).Value = 'S_Plant_growth_promoting'
$ws2.Cells.Item(25,3).Value = [double]'0.198578179801338'
$ws2.Cells.Item(25,4).Value = [double]'0.5'
$ws2.Cells.Item(25,5).Value = [double]'4'
$ws2.Cells.Item(26,1).Value = [double]'60'
$ws2.Cells.Item(26,2).Value = 'Unassigned'
$ws2.Cells.Item(26,3).Value = [double]'0.2619014764579588'
$ws2.Cells.Item(26,4).Value = [double]'0.2619047619047619'
$ws2.Cells.Item(26,5).Value = [double]'42'
$ws2.Cells.Item(27,1).Value = [double]'17'
$ws2.Cells.Item(27,2).Value = 'S_Transport__Inner_Membrane'
$ws2.Cells.Item(27,3).Value = [double]'0.2636771143363943'
$ws2.Cells.Item(27,4).Value = [double]'0.2304147465437788'
$ws2.Cells.Item(27,5).Value = [double]'217'
$ws2.Cells.Item(28,1).Value = [double]'10'
$ws2.Cells.Item(28,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Nicotinamide_Biosynthesis'
$ws2.Cells.Item(28,3).Value = [double]'0.2883750660594911'
$ws2.Cells.Item(28,4).Value = [double]'0.3076923076923077'
$ws2.Cells.Item(28,5).Value = [double]'13'
$ws2.Cells.Item(29,1).Value = [double]'3'
$ws2.Cells.Item(29,2).Value = 'S_Alternate_Carbon'
$ws2.Cells.Item(29,3).Value = [double]'0.2883750660594911'
$ws2.Cells.Item(29,4).Value = [double]'0.3076923076923077'
$ws2.Cells.Item(29,5).Value = [double]'13'
$ws2.Cells.Item(30,1).Value = [double]'88'
$ws2.Cells.Item(30,2).Value = 'S_tRNA_Charging'
$ws2.Cells.Item(30,3).Value = [double]'0.3560692766715254'
$ws2.Cells.Item(30,4).Value = [double]'0.2608695652173913'
$ws2.Cells.Item(30,5).Value = [double]'23'
$ws2.Cells.Item(31,1).Value = [double]'55'
$ws2.Cells.Item(31,2).Value = 'S_Formaldehyde_Metabolism'
$ws2.Cells.Item(31,3).Value = [double]'0.3728566694667126'
$ws2.Cells.Item(31,4).Value = [double]'0.3333333333333333'
$ws2.Cells.Item(31,5).Value = [double]'6'
$ws2.Cells.Item(32,1).Value = [double]'74'
$ws2.Cells.Item(32,2).Value = 'S_Pyruvate_Metabolism'
$ws2.Cells.Item(32,3).Value = [double]'0.3728566694667126'
$ws2.Cells.Item(32,4).Value = [double]'0.3333333333333333'
$ws2.Cells.Item(32,5).Value = [double]'6'
$ws2.Cells.Item(33,1).Value = [double]'29'
$ws2.Cells.Item(33,2).Value = 'S_Arginine_and_Proline_Metabolism'
$ws2.Cells.Item(33,3).Value = [double]'0.3784661058889651'
$ws2.Cells.Item(33,4).Value = [double]'0.25'
$ws2.Cells.Item(33,5).Value = [double]'28'
$ws2.Cells.Item(34,1).Value = [double]'20'
$ws2.Cells.Item(34,2).Value = 'S_Starch_and_Sucrose_Metabolism'
$ws2.Cells.Item(34,3).Value = [double]'0.3951078860410095'
$ws2.Cells.Item(34,4).Value = [double]'0.2666666666666667'
$ws2.Cells.Item(34,5).Value = [double]'15'
$ws2.Cells.Item(35,1).Value = [double]'32'
$ws2.Cells.Item(35,2).Value = 'S_Cell_Envelope_Biosynthesis__O_antigen_Biosynthesis'
$ws2.Cells.Item(35,3).Value = [double]'0.3951078860410095'
$ws2.Cells.Item(35,4).Value = [double]'0.2666666666666667'
$ws2.Cells.Item(35,5).Value = [double]'15'
$ws2.Cells.Item(36,1).Value = [double]'58'
$ws2.Cells.Item(36,2).Value = 'S_Urea_cycleamino_group_metabolism'
$ws2.Cells.Item(36,3).Value = [double]'0.4547569923788626'
$ws2.Cells.Item(36,4).Value = [double]'0.2857142857142857'
$ws2.Cells.Item(36,5).Value = [double]'7'
$ws2.Cells.Item(37,1).Value = [double]'44'
$ws2.Cells.Item(37,2).Value = 'S_Purine_Metabolism'
$ws2.Cells.Item(37,3).Value = [double]'0.4653727138953176'
$ws2.Cells.Item(37,4).Value = [double]'0.2207792207792208'
$ws2.Cells.Item(37,5).Value = [double]'77'
$ws2.Cells.Item(38,1).Value = [double]'61'
$ws2.Cells.Item(38,2).Value = 'S_Cysteine_Metabolism'
$ws2.Cells.Item(38,3).Value = [double]'0.5301174436155505'
$ws2.Cells.Item(38,4).Value = [double]'0.25'
$ws2.Cells.Item(38,5).Value = [double]'8'
$ws2.Cells.Item(39,1).Value = [double]'68'
$ws2.Cells.Item(39,2).Value = 'S_Glycolysis'
$ws2.Cells.Item(39,3).Value = [double]'0.5480319742498857'
$ws2.Cells.Item(39,4).Value = [double]'0.2222222222222222'
$ws2.Cells.Item(39,5).Value = [double]'18'
$ws2.Cells.Item(40,1).Value = [double]'7'
$ws2.Cells.Item(40,2).Value = 'S_Alternate_Carbon_and_Nitrogen_source__Dipeptide_Metabolism'
$ws2.Cells.Item(40,3).Value = [double]'0.5957208318225187'
$ws2.Cells.Item(40,4).Value = [double]'0.2083333333333333'
$ws2.Cells.Item(40,5).Value = [double]'24'
$ws2.Cells.Item(41,1).Value = [double]'21'
$ws2.Cells.Item(41,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis'
$ws2.Cells.Item(41,3).Value = [double]'0.613621147132835'
$ws2.Cells.Item(41,4).Value = [double]'0.25'
$ws2.Cells.Item(41,5).Value = [double]'4'
$ws2.Cells.Item(42,1).Value = [double]'70'
$ws2.Cells.Item(42,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Pantothenate_and_CoA_Biosynthesis'
$ws2.Cells.Item(42,3).Value = [double]'0.658293491333841'
$ws2.Cells.Item(42,4).Value = [double]'0.2'
$ws2.Cells.Item(42,5).Value = [double]'10'
$ws2.Cells.Item(43,1).Value = [double]'87'
$ws2.Cells.Item(43,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Vitamin_B6_Metabolism'
$ws2.Cells.Item(43,3).Value = [double]'0.658293491333841'
$ws2.Cells.Item(43,4).Value = [double]'0.2'
$ws2.Cells.Item(43,5).Value = [double]'10'
$ws2.Cells.Item(44,1).Value = [double]'9'
$ws2.Cells.Item(44,2).Value = 'S_Glyoxylate_and_dicarboxylate_metabolism'
$ws2.Cells.Item(44,3).Value = [double]'0.7110820023881432'
$ws2.Cells.Item(44,4).Value = [double]'0.1818181818181818'
$ws2.Cells.Item(44,5).Value = [double]'11'
$ws2.Cells.Item(45,1).Value = [double]'50'
$ws2.Cells.Item(45,2).Value = 'S_Xenobiotic_tolerance'
$ws2.Cells.Item(45,3).Value = [double]'0.7110820023881432'
$ws2.Cells.Item(45,4).Value = [double]'0.1818181818181818'
$ws2.Cells.Item(45,5).Value = [double]'11'
$ws2.Cells.Item(46,1).Value = [double]'67'
$ws2.Cells.Item(46,2).Value = 'S_Histidine_Metabolism'
$ws2.Cells.Item(46,3).Value = [double]'0.7305680843618261'
$ws2.Cells.Item(46,4).Value = [double]'0.1764705882352941'
$ws2.Cells.Item(46,5).Value = [double]'17'
$ws2.Cells.Item(47,1).Value = [double]'64'
$ws2.Cells.Item(47,2).Value = 'S_Nucleotide_Salvage_Pathway'
$ws2.Cells.Item(47,3).Value = [double]'0.7599615937812596'
$ws2.Cells.Item(47,4).Value = [double]'0.1666666666666667'
$ws2.Cells.Item(47,5).Value = [double]'6'
$ws2.Cells.Item(48,1).Value = [double]'56'
$ws2.Cells.Item(48,2).Value = 'S_Pyrimidine_Metabolism'
$ws2.Cells.Item(48,3).Value = [double]'0.7605707768410169'
$ws2.Cells.Item(48,4).Value = [double]'0.18'
$ws2.Cells.Item(48,5).Value = [double]'50'
$ws2.Cells.Item(49,1).Value = [double]'72'
$ws2.Cells.Item(49,2).Value = 'S_Cell_Envelope_Biosynthesis'
$ws2.Cells.Item(49,3).Value = [double]'0.7962344668455332'
$ws2.Cells.Item(49,4).Value = [double]'0.1538461538461539'
$ws2.Cells.Item(49,5).Value = [double]'13'
$ws2.Cells.Item(50,1).Value = [double]'45'
$ws2.Cells.Item(50,2).Value = 'S_Alternate_Carbon__Propanoate_Metabolism'
$ws2.Cells.Item(50,3).Value = [double]'0.8108290273483258'
$ws2.Cells.Item(50,4).Value = [double]'0.1428571428571428'
$ws2.Cells.Item(50,5).Value = [double]'7'
$ws2.Cells.Item(51,1).Value = [double]'90'
$ws2.Cells.Item(51,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__One_Carbon_pool_by_folate'
$ws2.Cells.Item(51,3).Value = [double]'0.8108290273483258'
$ws2.Cells.Item(51,4).Value = [double]'0.1428571428571428'
$ws2.Cells.Item(51,5).Value = [double]'7'
$ws2.Cells.Item(52,1).Value = [double]'39'
$ws2.Cells.Item(52,2).Value = 'S_Pentose_Phosphate_Pathway'
$ws2.Cells.Item(52,3).Value = [double]'0.8108290273483258'
$ws2.Cells.Item(52,4).Value = [double]'0.1428571428571428'
$ws2.Cells.Item(52,5).Value = [double]'7'
$ws2.Cells.Item(53,1).Value = [double]'41'
$ws2.Cells.Item(53,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Molybdopterin_Biosynthesis'
$ws2.Cells.Item(53,3).Value = [double]'0.8298559702662496'
$ws2.Cells.Item(53,4).Value = [double]'0.1428571428571428'
$ws2.Cells.Item(53,5).Value = [double]'14'
$ws2.Cells.Item(54,1).Value = [double]'23'
$ws2.Cells.Item(54,2).Value = 'S_Methionine_Metabolism'
$ws2.Cells.Item(54,3).Value = [double]'0.8298559702662496'
$ws2.Cells.Item(54,4).Value = [double]'0.1428571428571428'
$ws2.Cells.Item(54,5).Value = [double]'14'
$ws2.Cells.Item(55,1).Value = [double]'79'
$ws2.Cells.Item(55,2).Value = 'S_Phenylalanine_Tyrosine_Tryptophan_Biosynthesis'
$ws2.Cells.Item(55,3).Value = [double]'0.8540430221973266'
$ws2.Cells.Item(55,4).Value = [double]'0.1481481481481481'
$ws2.Cells.Item(55,5).Value = [double]'27'
$ws2.Cells.Item(56,1).Value = [double]'27'
$ws2.Cells.Item(56,2).Value = 'S_Glycine__Serine_and_threonine_metabolism'
$ws2.Cells.Item(56,3).Value = [double]'0.8583946062470124'
$ws2.Cells.Item(56,4).Value = [double]'0.1333333333333333'
$ws2.Cells.Item(56,5).Value = [double]'15'
$ws2.Cells.Item(57,1).Value = [double]'66'
$ws2.Cells.Item(57,2).Value = 'S_Glutamate_Metabolism'
$ws2.Cells.Item(57,3).Value = [double]'0.9027496339325705'
$ws2.Cells.Item(57,4).Value = [double]'0.1176470588235294'
$ws2.Cells.Item(57,5).Value = [double]'17'
$ws2.Cells.Item(58,1).Value = [double]'42'
$ws2.Cells.Item(58,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Queuosine_biosynthesis'
$ws2.Cells.Item(58,3).Value = [double]'0.9547669279587493'
$ws2.Cells.Item(58,4).Value = [double]'0.07692307692307693'
$ws2.Cells.Item(58,5).Value = [double]'13'
$ws2.Cells.Item(59,1).Value = [double]'26'
$ws2.Cells.Item(59,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Folate_Biosynthesis'
$ws2.Cells.Item(59,3).Value = [double]'0.9643754631658646'
$ws2.Cells.Item(59,4).Value = [double]'0.07142857142857142'
$ws2.Cells.Item(59,5).Value = [double]'14'
$ws2.Cells.Item(60,1).Value = [double]'49'
$ws2.Cells.Item(60,2).Value = 'S_Cell_Envelope_Biosynthesis__LPS_Biosynthesis'
$ws2.Cells.Item(60,3).Value = [double]'0.9643754631658646'
$ws2.Cells.Item(60,4).Value = [double]'0.07142857142857142'
$ws2.Cells.Item(60,5).Value = [double]'14'
$ws2.Cells.Item(61,1).Value = [double]'86'
$ws2.Cells.Item(61,2).Value = 'S_TCA_Cycle'
$ws2.Cells.Item(61,3).Value = [double]'0.9915147918716702'
$ws2.Cells.Item(61,4).Value = [double]'0.05'
$ws2.Cells.Item(61,5).Value = [double]'20'
$ws2.Cells.Item(62,1).Value = [double]'75'
$ws2.Cells.Item(62,2).Value = 'S_Aromatic_Compounds_Degradation__Toluene_Pathway'
$ws2.Cells.Item(62,3).Value = [double]'0.992948362645788'
$ws2.Cells.Item(62,4).Value = [double]'0.06666666666666667'
$ws2.Cells.Item(62,5).Value = [double]'30'
$ws2.Cells.Item(63,1).Value = [double]'5'
$ws2.Cells.Item(63,2).Value = 'S_Heavy_Metal_Tolerance'
$ws2.Cells.Item(63,3).Value = [double]'0.9953777531398144'
$ws2.Cells.Item(63,4).Value = [double]'0.0625'
$ws2.Cells.Item(63,5).Value = [double]'32'
$ws2.Cells.Item(64,1).Value = [double]'57'
$ws2.Cells.Item(64,2).Value = 'S_Iron_uptake_and_metabolism'
$ws2.Cells.Item(64,3).Value = [double]'0.9956114513724466'
$ws2.Cells.Item(64,4).Value = [double]'0.1052631578947368'
$ws2.Cells.Item(64,5).Value = [double]'76'
$ws2.Cells.Item(65,1).Value = [double]'35'
$ws2.Cells.Item(65,2).Value = 'S_Sulfur_Metabolism'
$ws2.Cells.Item(65,3).Value = [double]'0.9958640609350173'
$ws2.Cells.Item(65,4).Value = [double]'0.04347826086956522'
$ws2.Cells.Item(65,5).Value = [double]'23'
$ws2.Cells.Item(66,1).Value = [double]'62'
$ws2.Cells.Item(66,2).Value = 'S_Transport__Outer_Membrane'
$ws2.Cells.Item(66,3).Value = [double]'0.9998599633537646'
$ws2.Cells.Item(66,4).Value = [double]'0.1357615894039735'
$ws2.Cells.Item(66,5).Value = [double]'302'
$ws2.Cells.Item(67,1).Value = [double]'43'
$ws2.Cells.Item(67,2).Value = 'S_PHAs_Metabolism'
$ws2.Cells.Item(67,3).Value = [double]'0.9999999999989988'
$ws2.Cells.Item(67,4).Value = [double]'0.02631578947368421'
$ws2.Cells.Item(67,5).Value = [double]'152'
$ws2.Cells.Item(68,1).Value = [double]'6'
$ws2.Cells.Item(68,2).Value = 'Extracellular exchange'
$ws2.Cells.Item(68,3).Value = [double]'1'
$ws2.Cells.Item(68,4).Value = [double]'0'
$ws2.Cells.Item(68,5).Value = [double]'348'
$ws2.Cells.Item(69,1).Value = [double]'83'
$ws2.Cells.Item(69,2).Value = 'Aliphatic open-chain ketones metabolism'
$ws2.Cells.Item(69,3).Value = [double]'1'
$ws2.Cells.Item(69,4).Value = [double]'0'
$ws2.Cells.Item(69,5).Value = [double]'5'
$ws2.Cells.Item(70,1).Value = [double]'8'
$ws2.Cells.Item(70,2).Value = 'S_Butanoate_Metabolism'
$ws2.Cells.Item(70,3).Value = [double]'1'
$ws2.Cells.Item(70,4).Value = [double]'0'
$ws2.Cells.Item(70,5).Value = [double]'6'
$ws2.Cells.Item(71,1).Value = [double]'84'
$ws2.Cells.Item(71,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Ubiquinone_biosynthesis'
$ws2.Cells.Item(71,3).Value = [double]'1'
$ws2.Cells.Item(71,4).Value = [double]'0'
$ws2.Cells.Item(71,5).Value = [double]'10'
$ws2.Cells.Item(72,1).Value = [double]'11'
$ws2.Cells.Item(72,2).Value = 'S_Cell_Envelope_Biosynthesis__Biosynthesis_of_L_glycero_D_manno_heptose__Hep_'
$ws2.Cells.Item(72,3).Value = [double]'1'
$ws2.Cells.Item(72,4).Value = [double]'0'
$ws2.Cells.Item(72,5).Value = [double]'5'
$ws2.Cells.Item(73,1).Value = [double]'2'
$ws2.Cells.Item(73,2).Value = 'S_Aromatic_Compounds_Degradation__Gallic_acid_Metabolism'
$ws2.Cells.Item(73,3).Value = [double]'1'
$ws2.Cells.Item(73,4).Value = [double]'0'
$ws2.Cells.Item(73,5).Value = [double]'4'
$ws2.Cells.Item(74,1).Value = [double]'1'
$ws2.Cells.Item(74,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Thiamine_Metabolism'
$ws2.Cells.Item(74,3).Value = [double]'1'
$ws2.Cells.Item(74,4).Value = [double]'0'
$ws2.Cells.Item(74,5).Value = [double]'8'
$ws2.Cells.Item(75,1).Value = [double]'14'
$ws2.Cells.Item(75,2).Value = 'S_Oxidative_Phosphorylation'
$ws2.Cells.Item(75,3).Value = [double]'1'
$ws2.Cells.Item(75,4).Value = [double]'0'
$ws2.Cells.Item(75,5).Value = [double]'10'
$ws2.Cells.Item(76,1).Value = [double]'82'
$ws2.Cells.Item(76,2).Value = 'S_Cell_Envelope_Biosynthesis__Lipid_A_Biosynthesis'
$ws2.Cells.Item(76,3).Value = [double]'1'
$ws2.Cells.Item(76,4).Value = [double]'0'
$ws2.Cells.Item(76,5).Value = [double]'15'
$ws2.Cells.Item(77,1).Value = [double]'47'
$ws2.Cells.Item(77,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis___Pyrroloquinoline_quinone_biosynthesis'
$ws2.Cells.Item(77,3).Value = [double]'1'
$ws2.Cells.Item(77,4).Value = [double]'0'
$ws2.Cells.Item(77,5).Value = [double]'6'
$ws2.Cells.Item(78,1).Value = [double]'30'
$ws2.Cells.Item(78,2).Value = 'S_Alternate_Carbon__Ribose_Metabolism'
$ws2.Cells.Item(78,3).Value = [double]'1'
$ws2.Cells.Item(78,4).Value = [double]'0'
$ws2.Cells.Item(78,5).Value = [double]'5'
$ws2.Cells.Item(79,1).Value = [double]'73'
$ws2.Cells.Item(79,2).Value = 'S_Inorganic_polyphosphates_metabolism'
$ws2.Cells.Item(79,3).Value = [double]'1'
$ws2.Cells.Item(79,4).Value = [double]'0'
$ws2.Cells.Item(79,5).Value = [double]'7'
$ws2.Cells.Item(80,1).Value = [double]'69'
$ws2.Cells.Item(80,2).Value = 'S_Cell_Envelope_Biosynthesis__Cellulose_Metabolism'
$ws2.Cells.Item(80,3).Value = [double]'1'
$ws2.Cells.Item(80,4).Value = [double]'0'
$ws2.Cells.Item(80,5).Value = [double]'3'
$ws2.Cells.Item(81,1).Value = [double]'31'
$ws2.Cells.Item(81,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Biotin_Biosynthesis'
$ws2.Cells.Item(81,3).Value = [double]'1'
$ws2.Cells.Item(81,4).Value = [double]'0'
$ws2.Cells.Item(81,5).Value = [double]'15'
$ws2.Cells.Item(82,1).Value = [double]'34'
$ws2.Cells.Item(82,2).Value = 'S_Fatty_Acid__Biosynthesis'
$ws2.Cells.Item(82,3).Value = [double]'1'
$ws2.Cells.Item(82,4).Value = [double]'0'
$ws2.Cells.Item(82,5).Value = [double]'75'
$ws2.Cells.Item(83,1).Value = [double]'65'
$ws2.Cells.Item(83,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__FeS_clusters_metabolism'
$ws2.Cells.Item(83,3).Value = [double]'1'
$ws2.Cells.Item(83,4).Value = [double]'0'
$ws2.Cells.Item(83,5).Value = [double]'8'
$ws2.Cells.Item(84,1).Value = [double]'63'
$ws2.Cells.Item(84,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Riboflavin_Metabolism'
$ws2.Cells.Item(84,3).Value = [double]'1'
$ws2.Cells.Item(84,4).Value = [double]'0'
$ws2.Cells.Item(84,5).Value = [double]'9'
$ws2.Cells.Item(85,1).Value = [double]'37'
$ws2.Cells.Item(85,2).Value = 'S_Gluconeogenesis'
$ws2.Cells.Item(85,3).Value = [double]'1'
$ws2.Cells.Item(85,4).Value = [double]'0'
$ws2.Cells.Item(85,5).Value = [double]'12'
$ws2.Cells.Item(86,1).Value = [double]'59'
$ws2.Cells.Item(86,2).Value = 'S_Alanine_and_Aspartate_Metabolism'
$ws2.Cells.Item(86,3).Value = [double]'1'
$ws2.Cells.Item(86,4).Value = [double]'0'
$ws2.Cells.Item(86,5).Value = [double]'8'
$ws2.Cells.Item(87,1).Value = [double]'54'
$ws2.Cells.Item(87,2).Value = 'Biomass and maintenance functions'
$ws2.Cells.Item(87,3).Value = [double]'1'
$ws2.Cells.Item(87,4).Value = [double]'0'
$ws2.Cells.Item(87,5).Value = [double]'4'
$ws2.Cells.Item(88,1).Value = [double]'53'
$ws2.Cells.Item(88,2).Value = 'Intracellular demand'
$ws2.Cells.Item(88,3).Value = [double]'1'
$ws2.Cells.Item(88,4).Value = [double]'0'
$ws2.Cells.Item(88,5).Value = [double]'31'
$ws2.Cells.Item(89,1).Value = [double]'52'
$ws2.Cells.Item(89,2).Value = 'S_Cofactor_and_Prosthetic_Group_Biosynthesis__Terpenoid_backbone_biosynthesis'
$ws2.Cells.Item(89,3).Value = [double]'1'
$ws2.Cells.Item(89,4).Value = [double]'0'
$ws2.Cells.Item(89,5).Value = [double]'14'
$ws2.Cells.Item(90,1).Value = [double]'51'
$ws2.Cells.Item(90,2).Value = 'S_Phosphonate_and_phosphinate_metabolism'
$ws2.Cells.Item(90,3).Value = [double]'1'
$ws2.Cells.Item(90,4).Value = [double]'0'
$ws2.Cells.Item(90,5).Value = [double]'4'
$ws2.Cells.Item(91,1).Value = [double]'16'
$ws2.Cells.Item(91,2).Value = 'Intracellular source/sink'
$ws2.Cells.Item(91,3).Value = [double]'1'
$ws2.Cells.Item(91,4).Value = [double]'0'
$ws2.Cells.Item(91,5).Value = [double]'2'
$ws2.Cells.Item(92,1).Value = [double]'0'
$ws2.Cells.Item(92,2).Value = 'Murein Recycling'
$ws2.Cells.Item(92,3).Value = [double]'1'
$ws2.Cells.Item(92,4).Value = [double]'0'
$ws2.Cells.Item(92,5).Value = [double]'2'
